$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "PERFUMARIA" to "Sheet1"
$ws.Name = "Sheet1"

# Row 2 - dates/text in A, numeric values in B:G
$ws.Range("A2").Value = "'07/07/2023"
$ws.Range("A2").Style = "Normal"

$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 1000
$ws.Range("D2").Value = 1000
$ws.Range("E2").Value = 1000
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 100

# Row 3 - text values throughout (preserve as literal text, not auto-converted numbers/dates)
$ws.Range("A3").Value = "'07/07/2023"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").Value = "'2000.00"
$ws.Range("B3").Style = "Normal"

$ws.Range("C3").Value = "'3000.00"
$ws.Range("C3").Style = "Normal"

$ws.Range("D3").Value = "'2001.00"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = "'3001.00"
$ws.Range("E3").Style = "Normal"

$ws.Range("F3").Value = "'1.00"
$ws.Range("F3").Style = "Normal"

$ws.Range("G3").Value = "'100.03"
$ws.Range("G3").Style = "Normal"

# Remove row 4 entirely (shifts dimension from A1:G4 to A1:G3)
$ws.Rows(4).Delete()
